# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..G (TB, d2S, K, IP, Win, sum) for rows 2..13
$data = @{
    2  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
    3  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    4  = @(1.505614041169197, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 3.055818435266709)
    5  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    6  = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054)
    7  = @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 2.27892381503245)
    8  = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 9.576116808119359)
    9  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    10 = @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 1, 9.193893886484982)
    11 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126)
    12 = @(1.505614041169197, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 0, 2.495991366269901)
    13 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
